# Append: 2026-01-26 18:41 JST
# Update the "取得日時" (retrieved datetime) column (A) for all existing
# data rows on the "ランサーズ" sheet to the new timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-26 18:41:41"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value -ne $null -and "$($cell.Value)" -ne "") {
        $cell.Value = $newTimestamp
    }
}
